$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers I1 ("I0") and J1 ("IF"), matching the formatting of the
# existing header row (bold, bordered, centered) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-9 for the new I (I0) and J (IF) columns.
$iValues = @{2=1; 3=1; 4=1; 5=1; 6=1; 7=1; 8=8; 9=1}
$jValues = @{2=6; 3=5; 4=7; 5=5; 6=7; 7=5; 8=8; 9=1}

foreach ($r in 2..9) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
